$d = $word.ActiveDocument

$d.Content.Find.Execute("Astronomical Advancements Unveiling Cosmic Wonders", $false, $false, $false, $false, $false, $true, 1, $false, "Exploring the Realm of Science: Unveiling the Mysteries of Chemistry", 2) | Out-Null
$d.Content.Find.Execute("Dr. Evelyn Grant", $false, $false, $false, $false, $false, $true, 1, $false, "Clarissa Meyer", 2) | Out-Null
$d.Content.Find.Execute("grant@local", $false, $false, $false, $false, $false, $true, 1, $false, "clarissa.meyer@validschool", 2) | Out-Null
$d.Content.Find.Execute("From the ancient stargazers to modern-day astronomers, humanity's fascination with the cosmos has been an enduring thread throughout history", $false, $false, $false, $false, $false, $true, 1, $false, "Science, a captivating realm of discovery, unveils the enigmatic tapestry of our universe", 2) | Out-Null
$d.Content.Find.Execute(" As technology continues to advance, our understanding of the universe expands, revealing ever more profound mysteries and captivating discoveries", $false, $false, $false, $false, $false, $true, 1, $false, " As we traverse the vast expanse of knowledge, chemistry emerges as a cornerstone of understanding the intricate wonders of matter", 2) | Out-Null
$d.Content.Find.Execute(" This essay delves into the remarkable advancements in astronomy, highlighting how they have revolutionized our understanding of the vast expanse beyond our planet", $false, $false, $false, $false, $false, $true, 1, $false, " From the smallest atoms to the grand tapestry of chemical reactions, chemistry holds the key to unlocking nature's secrets. Embark on a journey through the captivating world of chemistry, where elements dance in harmony, molecules unfold their stories, and the symphony of chemical reactions unveils the mysteries of our existence", 2) | Out-Null
$d.Content.Find.Execute("Astronomers have pushed the boundaries of observation with the development of powerful telescopes and observatories, such as the Hubble Space Telescope and the Atacama Large Millimeter Array", $false, $false, $false, $false, $false, $true, 1, $false, "With an insatiable curiosity, chemists embark on a quest to comprehend the fundamental principles that govern the behavior of matter", 2) | Out-Null
$d.Content.Find.Execute(" These instruments have allowed us to peer deeper into space, uncovering distant galaxies, mysterious black holes, and breathtaking nebulas", $false, $false, $false, $false, $false, $true, 1, $false, " Like detectives unraveling an intricate puzzle, they meticulously analyze the properties and interactions of elements, revealing the intricate patterns and relationships hidden within the atomic realm", 2) | Out-Null
$d.Content.Find.Execute(" The insights gained from these observations have challenged our previous assumptions and opened up new avenues of exploration. This is truly an example of human ingenuity and scientific curiosity at its finest", $false, $false, $false, $false, $false, $true, 1, $false, " Through experiments, observations, and theoretical models, chemists strive to unveil the enigmatic choreography of atoms and molecules, unravelling the secrets of chemical reactions that shape our world", 2) | Out-Null
$d.Content.Find.Execute("Furthermore, the advent of space exploration missions has provided invaluable data and images, transforming our perception of celestial bodies within our solar system", $false, $false, $false, $false, $false, $true, 1, $false, "As we delve into the depths of chemistry, we encounter a world governed by precise laws and principles", 2) | Out-Null
$d.Content.Find.Execute(" Spacecraft like the Voyager probes and the Cassini-Huygens mission have sent back stunning visuals and data, revealing the intricacies of planets, moons, and asteroids", $false, $false, $false, $false, $false, $true, 1, $false, " From the periodic table, a symphony of elements arranged in a delicate dance of order, to the intricate ballet of electrons, protons, and neutrons within atoms, chemistry reveals the underlying elegance and harmony of the universe", 2) | Out-Null
$d.Content.Find.Execute(" These missions have not only expanded our knowledge of our cosmic neighborhood but also sparked a renewed sense of wonder and appreciation for the delicate balance of our own planet", $false, $false, $false, $false, $false, $true, 1, $false, " The study of chemistry empowers us to unravel the mysteries of matter's transformations, enabling us to harness the power of chemical reactions to create new materials, medicines, and technologies that shape our lives", 2) | Out-Null
$d.Content.Find.Execute("In conclusion, the advancements in astronomy have been nothing short of extraordinary, propelling us forward in our quest to understand the universe's intricate workings", $false, $false, $false, $false, $false, $true, 1, $false, "In this essay, we have embarked on an enchanting journey through the realm of chemistry, unveiling the enigmatic tapestry of matter's behavior", 2) | Out-Null
$d.Content.Find.Execute(" From the discovery of distant galaxies to the exploration of our own solar system, astronomers have pushed the boundaries of human knowledge and ignited our imaginations", $false, $false, $false, $false, $false, $true, 1, $false, " From the intricate dance of elements to the precise laws governing chemical reactions, chemistry emerges as a cornerstone of understanding the world around us", 2) | Out-Null
$d.Content.Find.Execute(" These advancements ", $false, $false, $false, $false, $false, $true, 1, $false, " Through experimentation and exploration, chemists unravel ", 2) | Out-Null
$d.Content.Find.Execute("serve as a testament to the indomitable spirit of exploration and the boundless possibilities that lie before us in the vast expanse of the cosmos", $false, $false, $false, $false, $false, $true, 1, $false, "the mysteries of atomic interactions, revealing the underlying harmony and elegance of the universe. As we continue to push the boundaries of chemical knowledge, we unlock the potential for transformative discoveries that will shape the future of science and technology, improving our lives and deepening our comprehension of the intricate wonders of our universe", 2) | Out-Null


# Append a new empty paragraph at the end of the document body
$sel = $word.Selection
$sel.EndKey(6) | Out-Null
$sel.TypeParagraph() | Out-Null

